# "contingencies with rene fine"
#
# lines_states.xlsx / Sheet1 holds a table of lines:
#   A=index  B=name  C=from_bus  D=to_bus  E=in_service
#
# Two new contingency lines ("line7", "line8") are inserted right after
# "line6" (pushing the "extr1".."extr8" rows down by two), and the
# from_bus / to_bus / in_service values for every affected row are
# refreshed to the new case data. The table grows from 14 to 16 data
# rows (A1:E15 -> A1:E17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rows 2-7 (index 0..5 / line1..line6) are untouched by this edit.
# rows 8-17 take on the new row layout:
#   row -> @(index, name, from_bus, to_bus, in_service)
$data = @(
  @(8,   6, "line7", 14, 11, $True),
  @(9,   7, "line8", 16,  9, $True),
  @(10,  8, "extr1",  5, 12, $True),
  @(11,  9, "extr2",  5,  9, $True),
  @(12, 10, "extr3", 10, 11, $True),
  @(13, 11, "extr4",  7,  8, $False),
  @(14, 12, "extr5",  9, 11, $True),
  @(15, 13, "extr6",  7, 11, $False),
  @(16, 14, "extr7",  5,  7, $True),
  @(17, 15, "extr8",  8,  5, $False)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
}

# rows 16-17 are brand new - give column A the same bold/bordered style
# used by the rest of the index column (A2:A15).
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
